# finish 3/5-3/10 data cleaning
# Collapse the two multi-line incident "Synopsis" descriptions (rows 19 & 20)
# into single-line text (replace embedded newline with a space), update the
# now-shorter row 20 to its new auto-fit height, and move the active
# selection from G20 to G19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (MO-20-S vandalism report): drop the internal line break.
$ws.Range("H19").Value = "(MO-20-S) REPORTS VANDALISM TO WINDOW VIA PUNCH, UNIT-3 DINING. SUSPECT: MALE, INDIAN, 18, 5'11, THIN, BLACK CURLY HAIR, WEARING GLASSES, PURPLE SWEATER, KHAKI PANTS, BLACK BACKPACK. LAST SEEN WALKING WESTBOUND FROM UNIT-3 DINING."

# Row 20 (SNEED, MALCOLM arrest report): drop the internal line break.
$ws.Range("H20").Value = "SNEED, MALCOLM (MB-22-O) ARRESTED FOR PROWLING, POSSESSION OF STOLEN PROPERTY AND VIOLATION OF PROBATION, BANWAY. TOT BPD JAIL. PROPERTY FROM CASE 19-00599 RECOVERED."

# Row 20 is now shorter (one fewer wrapped line) -- its row shrinks from 96 to 84pt.
$ws.Rows.Item(20).RowHeight = 84

# Move the active selection to G19 (was G20).
$ws.Range("G19").Select() | Out-Null
